# Remove form_id from remaining forms
#
# The "settings" sheet has a form_id column (column B) that is no longer
# wanted. Deleting it shifts version/style/namespaces left by one column
# (B<-C, C<-D, D<-E). Comments attached to the header cells are pinned to
# their literal cell address rather than following the data, so after the
# column delete we re-home each comment's text to the cell it now
# logically belongs to, and drop the now out-of-range trailing comment.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("settings")
$ws.Activate()

# Capture the comment text that needs to slide left before the column
# holding it disappears.
$versionComment    = $ws.Range("C1").Comment.Text()
$pagesComment       = $ws.Range("D1").Comment.Text()
$namespacesComment = $ws.Range("E1").Comment.Text()

# Delete the whole form_id column (B); C/D/E and their contents shift left.
$ws.Range("B:B").Delete()

# Re-home the comments onto the cells they now describe.
$ws.Range("B1").Comment.Text($versionComment)
$ws.Range("C1").Comment.Text($pagesComment)
$ws.Range("D1").Comment.Text($namespacesComment)
# The trailing comment (old E1 / namespaces) has nothing left to attach to.
$ws.Range("E1").Comment.Delete()

# Leave the settings sheet cursor on the new first data column.
$ws.Range("B1").Select()

# Restore focus to the survey sheet, with the frozen bottom-right pane's
# selection moved down to A2.
$ws1 = $wb.Worksheets.Item("survey")
$ws1.Activate()
$ws1.Range("A2").Select()
